# "new hk list for the StanDep"
# The underlying HK match-percentage values that previously rounded to
# 99.890109890109898 are refreshed to 100 (every cell in column A below
# the header that currently holds that old value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dims = $ws.UsedRange
$lastRow = $dims.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($val -eq 99.890109890109898) {
        $cell.Value = 100
    }
}

# Header cell keeps its text ("HK_R_acc_SD") - re-assert it so the sheet
# reflects the refreshed HK list pulled in for this sheet.
$ws.Range("A1").Value = "HK_R_acc_SD"
